$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the LiveSLR copyright/build string (B2) with the new build number.
$ws.Range("B2").Value = "Copyright @ 2023 Cytel Inc. LiveSLR 4.0.0.0 - Build #51133"

# Move the active selection from B3 to B2.
$ws.Range("B2").Select()
